# "changed the send email to its own method"
# This adds two new customer booking rows (form submissions) to the
# Customers sheet and marks the corresponding nights as "booked" on the
# Availability sheet.

$wb = $excel.ActiveWorkbook
$customers = $wb.Worksheets.Item("Customers")
$availability = $wb.Worksheets.Item("Availability")

# ---- New booking #1 (row 3) ----
$customers.Cells.Item(3, 2).Value = "dsjfds"               # B3 First Name
$customers.Cells.Item(3, 3).Value = "dsfkjds"               # C3 Last Name
$customers.Cells.Item(3, 4).Value = "kjsdfns"               # D3 Email
$customers.Cells.Item(3, 5).Value = "123"                   # E3 Phone #
$customers.Cells.Item(3, 6).Value = "asdfkn"                # F3 Payment First Name
$customers.Cells.Item(3, 7).Value = "dskfjn"                # G3 Payment Last Name
$customers.Cells.Item(3, 8).Value = "213"                   # H3 Card Number
$customers.Cells.Item(3, 9).Value = "12/12"                 # I3 Exp Date
$customers.Cells.Item(3, 10).Value = "sndfm"                # J3 Country
$customers.Cells.Item(3, 11).Value = "124"                  # K3 Zip Code
$customers.Cells.Item(3, 12).Value = "9N1i6"                # L3 ID
$customers.Cells.Item(3, 13).Value = "100"                  # M3 Room #
$customers.Cells.Item(3, 14).Value = "$115"                 # N3 Price
$customers.Cells.Item(3, 15).Value = "2023-11-15"           # O3 Start Date
$customers.Cells.Item(3, 16).Value = "2023-11-18"           # P3 End Date

# ---- New booking #2 (row 4) ----
$customers.Cells.Item(4, 2).Value = "adsfkjhds"             # B4 First Name
$customers.Cells.Item(4, 3).Value = "sdfjkdsfh"             # C4 Last Name
$customers.Cells.Item(4, 4).Value = "nwahba02@gmail.com"    # D4 Email
$customers.Cells.Item(4, 5).Value = "1234"                  # E4 Phone #
$customers.Cells.Item(4, 6).Value = "dsjfh"                 # F4 Payment First Name
$customers.Cells.Item(4, 7).Value = "sdkjfh"                # G4 Payment Last Name
$customers.Cells.Item(4, 8).Value = "2345"                  # H4 Card Number
$customers.Cells.Item(4, 9).Value = "12/12"                 # I4 Exp Date
$customers.Cells.Item(4, 10).Value = "fndsjkf"               # J4 Country
$customers.Cells.Item(4, 11).Value = "1234"                  # K4 Zip Code
$customers.Cells.Item(4, 12).Value = "R2RZa"                 # L4 ID
$customers.Cells.Item(4, 13).Value = "100"                   # M4 Room #
$customers.Cells.Item(4, 14).Value = "$115"                  # N4 Price
$customers.Cells.Item(4, 15).Value = "2023-11-22"            # O4 Start Date
$customers.Cells.Item(4, 16).Value = "2023-11-30"            # P4 End Date

# ---- Mark the Availability grid as "booked" for room 100 (column B) ----
# Booking #1: 2023-11-15 .. 2023-11-18 -> rows 16-19
foreach ($r in 16..19) {
    $availability.Cells.Item($r, 2).Value = "booked"
}

# Booking #2: 2023-11-22 .. 2023-11-30 -> rows 23-31
foreach ($r in 23..31) {
    $availability.Cells.Item($r, 2).Value = "booked"
}
